$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 298.29166
$ws.Range("I28").Value = 315.42105
$ws.Range("J28").Value = 233.2
$ws.Range("K28").Value = 315.42105
$ws.Range("L28").Value = 233.2
$ws.Range("M28").Value = 169.57895
$ws.Range("N28").Value = -1203.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2757.5881
$ws.Range("I131").Value = 2961.9092
$ws.Range("J131").Value = 2383
$ws.Range("K131").Value = 8885.7276
$ws.Range("L131").Value = 7149
$ws.Range("M131").Value = -3845.7276
$ws.Range("N131").Value = -17229

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3806.077
$ws.Range("I137").Value = 4136.6772
$ws.Range("K137").Value = 12410.0316
$ws.Range("M137").Value = -9860.0316

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1307.091
$ws.Range("I45").Value = 1048.6285
$ws.Range("J45").Value = 2312.2222
$ws.Range("K45").Value = 1048.6285
$ws.Range("L45").Value = 2312.2222
$ws.Range("M45").Value = -671.6285
$ws.Range("N45").Value = -3066.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2384.7
$ws.Range("I61").Value = 1076
$ws.Range("J61").Value = 3455.4546
$ws.Range("K61").Value = 1076
$ws.Range("L61").Value = 3455.4546
$ws.Range("M61").Value = -864
$ws.Range("N61").Value = -3879.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 546.9318
$ws.Range("I74").Value = 455.2439
$ws.Range("J74").Value = 1800
$ws.Range("K74").Value = 455.2439
$ws.Range("L74").Value = 1800
$ws.Range("M74").Value = 418.7561
$ws.Range("N74").Value = -3548

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 28684
$ws.Range("J76").Value = 28684
$ws.Range("L76").Value = 28684
$ws.Range("N76").Value = -29360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 546.9318
$ws.Range("I77").Value = 455.2439
$ws.Range("J77").Value = 1800
$ws.Range("K77").Value = 2276.2195
$ws.Range("L77").Value = 9000
$ws.Range("M77").Value = 2091.7805
$ws.Range("N77").Value = -17736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 28684
$ws.Range("J79").Value = 28684
$ws.Range("L79").Value = 28684
$ws.Range("N79").Value = -31024

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2237.4348
$ws.Range("I132").Value = 1507.2354
$ws.Range("J132").Value = 4306.3335
$ws.Range("K132").Value = 4521.706200000001
$ws.Range("L132").Value = 12919.0005
$ws.Range("M132").Value = -1991.706200000001
$ws.Range("N132").Value = -17979.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2384.7
$ws.Range("I136").Value = 1076
$ws.Range("J136").Value = 3455.4546
$ws.Range("K136").Value = 3228
$ws.Range("L136").Value = 10366.3638
$ws.Range("M136").Value = -678
$ws.Range("N136").Value = -15466.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 36714.168
$ws.Range("I76").Value = 20285
$ws.Range("K76").Value = 20285
$ws.Range("M76").Value = -19970

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H79").Value = 36714.168
$ws.Range("I79").Value = 20285
$ws.Range("K79").Value = 20285
$ws.Range("M79").Value = -19193

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4393.923
$ws.Range("I134").Value = 5629.0835
$ws.Range("J134").Value = 3335.2144
$ws.Range("K134").Value = 16887.2505
$ws.Range("L134").Value = 10005.6432
$ws.Range("M134").Value = -14352.2505
$ws.Range("N134").Value = -15075.6432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2513.5264
$ws.Range("I31").Value = 1479.3778
$ws.Range("J31").Value = 6391.5835
$ws.Range("K31").Value = 1479.3778
$ws.Range("L31").Value = 6391.5835
$ws.Range("M31").Value = -1184.3778
$ws.Range("N31").Value = -6981.5835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2513.5264
$ws.Range("I34").Value = 1479.3778
$ws.Range("J34").Value = 6391.5835
$ws.Range("K34").Value = 1479.3778
$ws.Range("L34").Value = 6391.5835
$ws.Range("M34").Value = -1277.3778
$ws.Range("N34").Value = -6795.5835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7355322
$ws.Range("I58").Value = 1529.0817
$ws.Range("K58").Value = 1529.0817
$ws.Range("M58").Value = -1326.0817

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 21006.8
$ws.Range("J74").Value = 21006.8
$ws.Range("L74").Value = 21006.8
$ws.Range("N74").Value = -22754.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 21006.8
$ws.Range("J77").Value = 21006.8
$ws.Range("L77").Value = 63020.39999999999
$ws.Range("N77").Value = -71756.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1586.3623
$ws.Range("I132").Value = 1268.86
$ws.Range("J132").Value = 2421.8948
$ws.Range("K132").Value = 3806.58
$ws.Range("L132").Value = 7265.6844
$ws.Range("M132").Value = -1276.58
$ws.Range("N132").Value = -12325.6844

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1574.7678
$ws.Range("I134").Value = 1069.9783
$ws.Range("J134").Value = 3896.8
$ws.Range("K134").Value = 3209.9349
$ws.Range("L134").Value = 11690.4
$ws.Range("M134").Value = -674.9349000000002
$ws.Range("N134").Value = -16760.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 7355322
$ws.Range("I136").Value = 1529.0817
$ws.Range("K136").Value = 4587.2451
$ws.Range("M136").Value = -2037.2451

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1892.2858
$ws.Range("J131").Value = 1525.9259
$ws.Range("L131").Value = 4577.7777
$ws.Range("N131").Value = -14657.7777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2613.9355
$ws.Range("I137").Value = 2294.8462
$ws.Range("J137").Value = 2844.389
$ws.Range("K137").Value = 6884.5386
$ws.Range("L137").Value = 8533.167000000001
$ws.Range("M137").Value = -1784.5386
$ws.Range("N137").Value = -18733.167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 34500
$ws.Range("J69").Value = 34500
$ws.Range("L69").Value = 34500
$ws.Range("N69").Value = -35998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H72").Value = 34500
$ws.Range("J72").Value = 34500
$ws.Range("L72").Value = 103500
$ws.Range("N72").Value = -110988

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 10311
$ws.Range("I123").Value = 10296
$ws.Range("J123").Value = 10326
$ws.Range("K123").Value = 10296
$ws.Range("L123").Value = 10326
$ws.Range("M123").Value = -7846
$ws.Range("N123").Value = -15226

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 17316.334
$ws.Range("J136").Value = 17316.334
$ws.Range("L136").Value = 51949.00199999999
$ws.Range("N136").Value = -57049.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 28661
$ws.Range("J76").Value = 28661
$ws.Range("L76").Value = 28661
$ws.Range("N76").Value = -29337

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 28661
$ws.Range("J79").Value = 28661
$ws.Range("L79").Value = 28661
$ws.Range("N79").Value = -31001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 36666.668
$ws.Range("J87").Value = 36666.668
$ws.Range("L87").Value = 36666.668
$ws.Range("N87").Value = -38912.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 26104.637
$ws.Range("I88").Value = 8575.5
$ws.Range("K88").Value = 8575.5
$ws.Range("M88").Value = -8147.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H90").Value = 36666.668
$ws.Range("J90").Value = 36666.668
$ws.Range("L90").Value = 110000.004
$ws.Range("N90").Value = -121232.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H91").Value = 26104.637
$ws.Range("I91").Value = 8575.5
$ws.Range("K91").Value = 8575.5
$ws.Range("M91").Value = -7093.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3007.0344
$ws.Range("I122").Value = 2800.2222
$ws.Range("J122").Value = 3345.4546
$ws.Range("K122").Value = 8400.6666
$ws.Range("L122").Value = 10036.3638
$ws.Range("M122").Value = -5950.6666
$ws.Range("N122").Value = -14936.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 28140
$ws.Range("J131").Value = 28140
$ws.Range("L131").Value = 28140
$ws.Range("N131").Value = -38220

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2332.756
$ws.Range("I132").Value = 1497.963
$ws.Range("J132").Value = 3942.7144
$ws.Range("K132").Value = 4493.889
$ws.Range("L132").Value = 11828.1432
$ws.Range("M132").Value = -1963.889
$ws.Range("N132").Value = -16888.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1804.8298
$ws.Range("I136").Value = 1290.081
$ws.Range("K136").Value = 3870.242999999999
$ws.Range("M136").Value = -1320.242999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2224999.8
$ws.Range("I3").Value = 5500000
$ws.Range("J3").Value = 41666.332
$ws.Range("K3").Value = 5500000
$ws.Range("L3").Value = 41666.332
$ws.Range("M3").Value = -5499886
$ws.Range("N3").Value = -41894.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7431.9316
$ws.Range("I132").Value = 1674.193
$ws.Range("J132").Value = 18018.742
$ws.Range("K132").Value = 5022.579
$ws.Range("L132").Value = 54056.226
$ws.Range("M132").Value = -2492.579
$ws.Range("N132").Value = -59116.226

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 766.1111
$ws.Range("I136").Value = 613.6957
$ws.Range("J136").Value = 1035.7693
$ws.Range("K136").Value = 1841.0871
$ws.Range("L136").Value = 3107.3079
$ws.Range("M136").Value = 708.9129
$ws.Range("N136").Value = -8207.3079
